# Sync the ETF ticker list:
#  - Remove the unavailable KOTAKLOVOL row entirely
#  - Rename the remaining Kotak-prefixed tickers to their new short names
#    (KOTAKALPHA -> ALPHA, KOTAKIT -> IT, KOTAKPSUBK -> PSUBANK, KOTAKNV20 -> NV20)
#  - Leave the selection on the row that used to hold KOTAKLOVOL, matching
#    the cursor position left behind by deleting that row in the UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# KOTAKLOVOL has no replacement ETF any more - drop its whole row so every
# row beneath it shifts up.
$lovolCell = $ws.Cells.Find("KOTAKLOVOL")
$lovolRow = $lovolCell.Row
$lovolCell.EntireRow.Delete()

# Rename the remaining Kotak ETFs to their new tickers, preserving each
# row's historical data / formatting.
$ws.Cells.Find("KOTAKALPHA").Value = "ALPHA"
$ws.Cells.Find("KOTAKIT").Value = "IT"
$ws.Cells.Find("KOTAKPSUBK").Value = "PSUBANK"
$ws.Cells.Find("KOTAKNV20").Value = "NV20"

# Leave the selection where the deleted row used to be (whole row selected),
# matching the post-delete cursor position.
[void]$ws.Rows($lovolRow).Select()
